$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-18 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-19 Friday", 2) | Out-Null
$d.Content.Find.Execute("36+48=", $true, $false, $false, $false, $false, $true, 1, $false, "58-4=", 2) | Out-Null
$d.Content.Find.Execute("97-52=", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=", 2) | Out-Null
$d.Content.Find.Execute("25+15=", $true, $false, $false, $false, $false, $true, 1, $false, "73-17=", 2) | Out-Null
$d.Content.Find.Execute("84-42=", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=", 2) | Out-Null
$d.Content.Find.Execute("94-85=", $true, $false, $false, $false, $false, $true, 1, $false, "7+15=", 2) | Out-Null
$d.Content.Find.Execute("65-58=", $true, $false, $false, $false, $false, $true, 1, $false, "33-30=", 2) | Out-Null
$d.Content.Find.Execute("49-9=", $true, $false, $false, $false, $false, $true, 1, $false, "97-84=", 2) | Out-Null
$d.Content.Find.Execute("41-24=", $true, $false, $false, $false, $false, $true, 1, $false, "2+79=", 2) | Out-Null
$d.Content.Find.Execute("41+56=", $true, $false, $false, $false, $false, $true, 1, $false, "42-7=", 2) | Out-Null
$d.Content.Find.Execute("11+62=", $true, $false, $false, $false, $false, $true, 1, $false, "88-15=", 2) | Out-Null
$d.Content.Find.Execute("40+0=", $true, $false, $false, $false, $false, $true, 1, $false, "66-52=", 2) | Out-Null
$d.Content.Find.Execute("26-18=", $true, $false, $false, $false, $false, $true, 1, $false, "51-2=", 2) | Out-Null
$d.Content.Find.Execute("70-20=", $true, $false, $false, $false, $false, $true, 1, $false, "0+44=", 2) | Out-Null
$d.Content.Find.Execute("90-40=", $true, $false, $false, $false, $false, $true, 1, $false, "16+44=", 2) | Out-Null
$d.Content.Find.Execute("6-6=", $true, $false, $false, $false, $false, $true, 1, $false, "85-12=", 2) | Out-Null
$d.Content.Find.Execute("42-19=", $true, $false, $false, $false, $false, $true, 1, $false, "6-4=", 2) | Out-Null
$d.Content.Find.Execute("14+5=", $true, $false, $false, $false, $false, $true, 1, $false, "4+14=", 2) | Out-Null
$d.Content.Find.Execute("10+54=", $true, $false, $false, $false, $false, $true, 1, $false, "76-50=", 2) | Out-Null
$d.Content.Find.Execute("3+33=", $true, $false, $false, $false, $false, $true, 1, $false, "36-13=", 2) | Out-Null
$d.Content.Find.Execute("66+28=", $true, $false, $false, $false, $false, $true, 1, $false, "52-11=", 2) | Out-Null
$d.Content.Find.Execute("78-58=", $true, $false, $false, $false, $false, $true, 1, $false, "63-60=", 2) | Out-Null
$d.Content.Find.Execute("51-24=", $true, $false, $false, $false, $false, $true, 1, $false, "70-10=", 2) | Out-Null
$d.Content.Find.Execute("35+63=", $true, $false, $false, $false, $false, $true, 1, $false, "41+6=", 2) | Out-Null
$d.Content.Find.Execute("41-8=", $true, $false, $false, $false, $false, $true, 1, $false, "49-37=", 2) | Out-Null
$d.Content.Find.Execute("20-19=", $true, $false, $false, $false, $false, $true, 1, $false, "22+38=", 2) | Out-Null
$d.Content.Find.Execute("73-6=", $true, $false, $false, $false, $false, $true, 1, $false, "1+3=", 2) | Out-Null
$d.Content.Find.Execute("77-20=", $true, $false, $false, $false, $false, $true, 1, $false, "47+3=", 2) | Out-Null
$d.Content.Find.Execute("95-62=", $true, $false, $false, $false, $false, $true, 1, $false, "76+11=", 2) | Out-Null
$d.Content.Find.Execute("61-6=", $true, $false, $false, $false, $false, $true, 1, $false, "70-28=", 2) | Out-Null
$d.Content.Find.Execute("72-43=", $true, $false, $false, $false, $false, $true, 1, $false, "8+49=", 2) | Out-Null
$d.Content.Find.Execute("53+38=", $true, $false, $false, $false, $false, $true, 1, $false, "80-40=", 2) | Out-Null
$d.Content.Find.Execute("77-41=", $true, $false, $false, $false, $false, $true, 1, $false, "52-24=", 2) | Out-Null
$d.Content.Find.Execute("6+79=", $true, $false, $false, $false, $false, $true, 1, $false, "20+51=", 2) | Out-Null
$d.Content.Find.Execute("40+1=", $true, $false, $false, $false, $false, $true, 1, $false, "88-53=", 2) | Out-Null
$d.Content.Find.Execute("76-37=", $true, $false, $false, $false, $false, $true, 1, $false, "71+9=", 2) | Out-Null
$d.Content.Find.Execute("41-7=", $true, $false, $false, $false, $false, $true, 1, $false, "66+5=", 2) | Out-Null
$d.Content.Find.Execute("22-15=", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=", 2) | Out-Null
$d.Content.Find.Execute("10+8=", $true, $false, $false, $false, $false, $true, 1, $false, "76-21=", 2) | Out-Null
$d.Content.Find.Execute("50-11=", $true, $false, $false, $false, $false, $true, 1, $false, "74-17=", 2) | Out-Null
$d.Content.Find.Execute("0+36=", $true, $false, $false, $false, $false, $true, 1, $false, "27+12=", 2) | Out-Null
$d.Content.Find.Execute("21+77=", $true, $false, $false, $false, $false, $true, 1, $false, "61-55=", 2) | Out-Null
$d.Content.Find.Execute("74-14=", $true, $false, $false, $false, $false, $true, 1, $false, "15+1=", 2) | Out-Null
$d.Content.Find.Execute("42+38=", $true, $false, $false, $false, $false, $true, 1, $false, "29+45=", 2) | Out-Null
$d.Content.Find.Execute("91-62=", $true, $false, $false, $false, $false, $true, 1, $false, "28+46=", 2) | Out-Null
$d.Content.Find.Execute("25+61=", $true, $false, $false, $false, $false, $true, 1, $false, "80-80=", 2) | Out-Null
$d.Content.Find.Execute("99-58=", $true, $false, $false, $false, $false, $true, 1, $false, "20+76=", 2) | Out-Null
$d.Content.Find.Execute("8+91=", $true, $false, $false, $false, $false, $true, 1, $false, "24+49=", 2) | Out-Null
$d.Content.Find.Execute("78-28=", $true, $false, $false, $false, $false, $true, 1, $false, "62-9=", 2) | Out-Null
$d.Content.Find.Execute("60-36=", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=", 2) | Out-Null
$d.Content.Find.Execute("22+16=", $true, $false, $false, $false, $false, $true, 1, $false, "25+24=", 2) | Out-Null
$d.Content.Find.Execute("25-19=", $true, $false, $false, $false, $false, $true, 1, $false, "27+21=", 2) | Out-Null
$d.Content.Find.Execute("50+46=", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=", 2) | Out-Null
$d.Content.Find.Execute("69+7=", $true, $false, $false, $false, $false, $true, 1, $false, "86-10=", 2) | Out-Null
$d.Content.Find.Execute("27+39=", $true, $false, $false, $false, $false, $true, 1, $false, "65-55=", 2) | Out-Null
$d.Content.Find.Execute("68-18=", $true, $false, $false, $false, $false, $true, 1, $false, "73+21=", 2) | Out-Null
$d.Content.Find.Execute("34+46=", $true, $false, $false, $false, $false, $true, 1, $false, "50-36=", 2) | Out-Null
$d.Content.Find.Execute("88+0=", $true, $false, $false, $false, $false, $true, 1, $false, "79-69=", 2) | Out-Null
$d.Content.Find.Execute("4+95=", $true, $false, $false, $false, $false, $true, 1, $false, "92-48=", 2) | Out-Null
$d.Content.Find.Execute("58-10=", $true, $false, $false, $false, $false, $true, 1, $false, "28+33=", 2) | Out-Null
$d.Content.Find.Execute("84-83=", $true, $false, $false, $false, $false, $true, 1, $false, "43+40=", 2) | Out-Null
$d.Content.Find.Execute("59-47=", $true, $false, $false, $false, $false, $true, 1, $false, "69-13=", 2) | Out-Null
$d.Content.Find.Execute("84-78=", $true, $false, $false, $false, $false, $true, 1, $false, "69-15=", 2) | Out-Null
$d.Content.Find.Execute("31-6=", $true, $false, $false, $false, $false, $true, 1, $false, "24+8=", 2) | Out-Null
$d.Content.Find.Execute("41-13=", $true, $false, $false, $false, $false, $true, 1, $false, "3+86=", 2) | Out-Null
$d.Content.Find.Execute("53-20=", $true, $false, $false, $false, $false, $true, 1, $false, "22-12=", 2) | Out-Null
$d.Content.Find.Execute("29+6=", $true, $false, $false, $false, $false, $true, 1, $false, "90-46=", 2) | Out-Null
$d.Content.Find.Execute("96-14=", $true, $false, $false, $false, $false, $true, 1, $false, "51-17=", 2) | Out-Null
$d.Content.Find.Execute("48-2=", $true, $false, $false, $false, $false, $true, 1, $false, "18+31=", 2) | Out-Null
$d.Content.Find.Execute("3+64=", $true, $false, $false, $false, $false, $true, 1, $false, "5+89=", 2) | Out-Null
$d.Content.Find.Execute("45+0=", $true, $false, $false, $false, $false, $true, 1, $false, "85+6=", 2) | Out-Null
$d.Content.Find.Execute("29-18=", $true, $false, $false, $false, $false, $true, 1, $false, "77+20=", 2) | Out-Null
$d.Content.Find.Execute("4+38=", $true, $false, $false, $false, $false, $true, 1, $false, "68-12=", 2) | Out-Null
$d.Content.Find.Execute("86-43=", $true, $false, $false, $false, $false, $true, 1, $false, "76-53=", 2) | Out-Null
$d.Content.Find.Execute("58+21=", $true, $false, $false, $false, $false, $true, 1, $false, "31+68=", 2) | Out-Null
$d.Content.Find.Execute("86-51=", $true, $false, $false, $false, $false, $true, 1, $false, "47-33=", 2) | Out-Null
$d.Content.Find.Execute("30-11=", $true, $false, $false, $false, $false, $true, 1, $false, "71-5=", 2) | Out-Null
$d.Content.Find.Execute("31+15=", $true, $false, $false, $false, $false, $true, 1, $false, "60-60=", 2) | Out-Null
$d.Content.Find.Execute("11-8=", $true, $false, $false, $false, $false, $true, 1, $false, "35+20=", 2) | Out-Null
$d.Content.Find.Execute("7+33=", $true, $false, $false, $false, $false, $true, 1, $false, "89-57=", 2) | Out-Null
$d.Content.Find.Execute("97-25=", $true, $false, $false, $false, $false, $true, 1, $false, "57-52=", 2) | Out-Null
$d.Content.Find.Execute("38+30=", $true, $false, $false, $false, $false, $true, 1, $false, "59-38=", 2) | Out-Null
$d.Content.Find.Execute("84-22=", $true, $false, $false, $false, $false, $true, 1, $false, "63+31=", 2) | Out-Null
$d.Content.Find.Execute("37-31=", $true, $false, $false, $false, $false, $true, 1, $false, "42-16=", 2) | Out-Null
$d.Content.Find.Execute("67-53=", $true, $false, $false, $false, $false, $true, 1, $false, "0+74=", 2) | Out-Null
$d.Content.Find.Execute("70-32=", $true, $false, $false, $false, $false, $true, 1, $false, "38-6=", 2) | Out-Null
$d.Content.Find.Execute("50+14=", $true, $false, $false, $false, $false, $true, 1, $false, "43+28=", 2) | Out-Null
$d.Content.Find.Execute("49-28=", $true, $false, $false, $false, $false, $true, 1, $false, "99-37=", 2) | Out-Null
$d.Content.Find.Execute("7+77=", $true, $false, $false, $false, $false, $true, 1, $false, "71-67=", 2) | Out-Null
$d.Content.Find.Execute("98-52=", $true, $false, $false, $false, $false, $true, 1, $false, "96-22=", 2) | Out-Null
$d.Content.Find.Execute("95-54=", $true, $false, $false, $false, $false, $true, 1, $false, "45-41=", 2) | Out-Null
$d.Content.Find.Execute("14+84=", $true, $false, $false, $false, $false, $true, 1, $false, "93-39=", 2) | Out-Null
$d.Content.Find.Execute("67+28=", $true, $false, $false, $false, $false, $true, 1, $false, "72+22=", 2) | Out-Null
$d.Content.Find.Execute("12-5=", $true, $false, $false, $false, $false, $true, 1, $false, "86-57=", 2) | Out-Null
$d.Content.Find.Execute("66-59=", $true, $false, $false, $false, $false, $true, 1, $false, "32-12=", 2) | Out-Null
$d.Content.Find.Execute("96+0=", $true, $false, $false, $false, $false, $true, 1, $false, "98-71=", 2) | Out-Null
$d.Content.Find.Execute("22+25=", $true, $false, $false, $false, $false, $true, 1, $false, "3+96=", 2) | Out-Null
$d.Content.Find.Execute("92-27=", $true, $false, $false, $false, $false, $true, 1, $false, "63+11=", 2) | Out-Null
$d.Content.Find.Execute("48+30=", $true, $false, $false, $false, $false, $true, 1, $false, "18+4=", 2) | Out-Null
$d.Content.Find.Execute("54+39=", $true, $false, $false, $false, $false, $true, 1, $false, "0+82=", 2) | Out-Null
$d.Content.Find.Execute("18+79=", $true, $false, $false, $false, $false, $true, 1, $false, "25-13=", 2) | Out-Null
